$d = $word.ActiveDocument

# The document currently opens with a single bold paragraph reading
# "剧情框架（主线：60分钟倒计时）". The edit prepends a bold file-tree listing
# (a "/js" folder with config.js / utils.js / minigames.js / engine.js /
# main.js plus a trailing blank line) as seven new paragraphs in front of
# that heading, leaving the heading paragraph itself untouched.

$first = $d.Paragraphs.Item(1)
$insertionPoint = $first.Range
$insertionPoint.Collapse(1)

$startPos = $insertionPoint.Start

$treeLines = @(
    "/js",
    " ├─ config.js       // 游戏配置 & 路线数据",
    " ├─ utils.js        // 工具函数（mmss、logMsg、spend、setState、renderClock…）",
    " ├─ minigames.js    // 小游戏统一入口 startMiniGame(...)",
    " ├─ engine.js       // 流程引擎（startCountdown、chooseRoute、renderStep、finishJourney…）",
    " └─ main.js         // 入口文件（事件绑定 & 初始化）",
    ""
)
$treeText = ($treeLines -join "`r") + "`r"

$insertionPoint.InsertBefore($treeText)

$endPos = $startPos + $treeText.Length
$newRange = $d.Range($startPos, $endPos)
$newRange.Font.Bold = 1
$newRange.Font.BoldBi = 1

Write-Output "paragraphs now: $($d.Paragraphs.Count)"
